$d = $word.ActiveDocument

# Locate the paragraph that currently holds only the "_GoBack" bookmark
# (an otherwise-empty paragraph right after the "Presentacion Personal"
# bio paragraph that ends in "...Soy Esposa, Mama y Abuela..."). We
# replace it with a blank paragraph followed by the new self-
# introduction paragraphs, then a trailing blank paragraph, exactly as
# in the target revision.
$idx = 0
$bioIndex = -1
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -match "Abuela") {
        $bioIndex = $idx
    }
}
if ($bioIndex -gt 0) {
    $target = $d.Paragraphs($bioIndex + 1)
} else {
    # Fixed fallback matching the known document layout.
    $target = $d.Paragraphs(14)
}

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:lang w:val="es-AR"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="es-AR"/>
    </w:rPr>
    <w:t xml:space="preserve">Hoy quiero hablarles un poco sobre m&#237;, mi trabajo y lo que creo. Soy </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="es-AR"/>
    </w:rPr>
    <w:t xml:space="preserve">Luis Salazar </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="es-AR"/>
    </w:rPr>
    <w:t>programador y actualmente trabajo en una empresa privada donde tengo la oportunidad de enfrentar desaf&#237;os diariamente. Me encanta el hecho de poder encontrar soluciones creativas y efectivas para los problemas que se me presentan. Creo firmemente que la programaci&#243;n es una herramienta poderosa que puede cambiar el mundo y estoy emocionado de poder ser parte de ese cambio.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:lang w:val="es-AR"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="es-AR"/>
    </w:rPr>
    <w:t>En mi trabajo, me esfuerzo por estar a la vanguardia de las &#250;ltimas tecnolog&#237;as y tendencias en programaci&#243;n. Siempre busco formas de mejorar mi conocimiento y habilidades para poder ser m&#225;s eficiente y efectivo en mi trabajo. Soy una persona que disfruta del aprendizaje continuo y estoy convencido de que nunca se termina de aprender, especialmente en el campo de la tecnolog&#237;a donde todo evoluciona muy r&#225;pido.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:lang w:val="es-AR"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="es-AR"/>
    </w:rPr>
    <w:t>Adem&#225;s, me gustan los desaf&#237;os. Me emociona enfrentar problemas complicados y poder encontrar soluciones innovadoras para ellos. Siempre trato de pensar fuera de la caja y encontrar una forma diferente de resolver los problemas.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:lang w:val="es-AR"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="es-AR"/>
    </w:rPr>
    <w:t>Por &#250;ltimo, quiero mencionar que no creo en los l&#237;mites cognitivos. Creo que todos tenemos la capacidad de aprender y ser buenos en lo que nos proponemos. Es cuesti&#243;n de pr&#225;ctica, dedicaci&#243;n y esfuerzo. Siempre estoy buscando formas de mejorar y superar mis propios l&#237;mites.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:lang w:val="es-AR"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="es-AR"/>
    </w:rPr>
    <w:t xml:space="preserve">En </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="es-AR"/>
    </w:rPr>
    <w:t>resumen</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="es-AR"/>
    </w:rPr>
    <w:t xml:space="preserve"> me encantan los desaf&#237;os y no creo en los l&#237;mites cognitivos. Estoy emocionado de formar parte de una comunidad de personas que comparten mi pasi&#243;n y entusiasmo por la programaci&#243;n</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="es-AR"/>
    </w:rPr>
    <w:t>.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:vanish/>
      <w:lang w:val="es-AR"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:vanish/>
      <w:lang w:val="es-AR"/>
    </w:rPr>
    <w:t>Principio del formulario</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@

$target.Range.InsertXML($xml)
